$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Repayment schedule": insert a new populated column O
# (mirrors column N's values/format for rows 2-15) and update the
# view selection to the full row 16.
# ---------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Copy the formatting of column N (rows 2-15) onto column O so the
# new cells pick up the same style (vertical-center + wrap text).
$wsRepay.Range("N2:N15").Copy()
$wsRepay.Range("O2:O15").PasteSpecial(-4122)
$wsRepay.Application.CutCopyMode = $false

# Populate the values that mirror column N (0 for data rows, blank
# for the "disbursement" spacer rows 2 and 4).
$wsRepay.Range("O3").Value2 = 0
$wsRepay.Range("O5:O15").Value2 = 0

# Move the selection to the full 16th row.
$wsRepay.Range("A16:XFD16").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet "Transactions": renumber the transaction IDs and move the
# view selection to D2.
# ---------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")

$wsTxn.Range("A2").Value2 = 76
$wsTxn.Range("A3").Value2 = 75
$wsTxn.Range("A4").Value2 = 74

$wsTxn.Range("D2").Select() | Out-Null
